# Insert a new data row at row 537 (pushing the existing rows 537-640
# down to 538-641) and populate it with the new weekly price record for
# Brócoli / Macroferia Regional de Talca.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 537..640 down by one, copying formatting (date style, etc.)
# from the row above - this is exactly what Excel's own "Insert" does.
$ws.Rows.Item(537).Insert()

# Populate the freshly inserted row with the new record's data.
$ws.Range("A537").Value = 5
$ws.Range("B537").Value = "Macroferia Regional de Talca"
$ws.Range("C537").Value = "Maule"
$ws.Range("D537").Value = "2023-10-10"
$ws.Range("E537").Value = 7
$ws.Range("F537").Value = 100112023
$ws.Range("G537").Value = "Brócoli"
$ws.Range("H537").Value = "Sin especificar"
$ws.Range("I537").Value = "Primera"
$ws.Range("J537").Value = 3000
$ws.Range("K537").Value = 1000
$ws.Range("L537").Value = 1000
$ws.Range("M537").Value = 1000
$ws.Range("N537").Value = "$/unidad"
$ws.Range("O537").Value = "Región del Maule"
$ws.Range("P537").Value = 1000
$ws.Range("Q537").Value = 1
$ws.Range("R537").Value = "Hortaliza"
